$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row at row 42 (shifts everything below it, including
# the September_Details/September_Date entries in columns R/S and the
# "Broadband" label further down in column A, down by one row).
$ws.Rows("42:42").Insert()

# Populate the newly inserted row with the new September entry.
$ws.Range("R42").Value = "balance your axis"
$ws.Range("S42").Value = "2024-09-18 10:28:28"
